$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D4").Value = "additional display,stock"
$ws.Range("D5").Select()
